$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CPF value looks numeric but must be stored as text (matching how the
# other CPF values in column G are stored), so force a text format before
# assigning it. Doing this before the style/format copy below lets the
# later paste normalize the cell style back in line with the rest of the
# row.
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Range("G7").Value = "12345678915"

# Copy the formatting of the last existing data row (row 6) down into the
# new row 7 so the new row matches the styling of the other data rows.
$ws.Range("A6:H6").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the new user's data.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Joao Pedro"
$ws.Range("C7").Value = "Santos Costa"
$ws.Range("D7").Value = "gfgn@dingf.com"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Masculino"
$ws.Range("H7").Value = "Aa123456789*"
